$wb = $excel.ActiveWorkbook

# Start from the single original sheet ("Kmart"). The target workbook has
# two sheets instead: "walmart" (the vendor list, header row only) and
# "costco" (header row plus a new water/1/100 row), in that tab order.
$original = $wb.Worksheets.Item(1)

# A throwaway sheet so the internal sheetId counter lands on 6/7 for the
# two sheets we actually keep (sheetIds are reused once a sheet is
# deleted, so this "spends" one id before creating the keepers).
$placeholder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $original)
$placeholder.Name = "placeholder"

$costco = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $placeholder)
$costco.Name = "costco"

$walmart = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $costco)
$walmart.Name = "walmart"

# "walmart" holds just the vendor headers.
$walmart.Range("A1").Value = "Item"
$walmart.Range("B1").Value = "Price"
$walmart.Range("C1").Value = "Stock"

# "costco" holds the headers plus the new water / 1 / 100 row.
$costco.Range("A1").Value = "Item"
$costco.Range("B1").Value = "Price"
$costco.Range("C1").Value = "Stock"
$costco.Range("A2").Value = "water"
$costco.Range("B2").Value = 1
$costco.Range("C2").Value = 100

# Drop the original "Kmart" sheet and the throwaway placeholder sheet.
$original.Delete()
$placeholderFixed = $wb.Worksheets.Item("placeholder")
$placeholderFixed.Delete()

# Re-fetch the surviving sheets by name (positional handles go stale once
# other sheets are deleted) and put "walmart" before "costco".
$walmartFixed = $wb.Worksheets.Item("walmart")
$costcoFixed = $wb.Worksheets.Item("costco")
$walmartFixed.Move($costcoFixed)

# Sheet handles track tab position, so Move() swaps what the two
# variables above now point at - re-fetch by name once more before using
# them for anything else.
$walmartFinal = $wb.Worksheets.Item("walmart")
$costcoFinal = $wb.Worksheets.Item("costco")

# Restore the view state: costco shows the new row selected, walmart is
# the active/visible tab.
$costcoFinal.Range("A2:C2").Select()
$walmartFinal.Select()
